# Commit: Add countries lists for the Headers organization and translations
# Fix the wording of the "Waktu Kampanye" (Campaign Time) line that appears
# throughout the document: move "2022" next to "Kampanye" and insert "untuk"
# before "rasi bintang gemini".

$d = $word.ActiveDocument

$find = "Waktu Kampanye rasi bintang gemini 2022: 14-23 Februari, 14-24 Maret"
$replace = "Waktu Kampanye 2022 untuk rasi bintang gemini: 14-23 Februari, 14-24 Maret"

# wdFindContinue=1, wdReplaceAll=2 -- replaces every occurrence in the document body.
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
